# Show dynamic msg to UI
# - Typography sheet: set "Wildcard Ranges" for the Default typography (row 4) to 0x20-0x7E
# - Translation sheet: add two new translation rows (SingleUseId2 / SingleUseId3)
#   carrying the new "RX: <value>" and "default" UI strings.

$wb = $excel.ActiveWorkbook

$wsTypography = $wb.Worksheets.Item("Typography")
$wsTypography.Range("I4").Value = "0x20-0x7E"

$wsTranslation = $wb.Worksheets.Item("Translation")

$wsTranslation.Range("B5").Value = "SingleUseId2"
$wsTranslation.Range("C5").Value = "Default"
$wsTranslation.Range("D5").Value = "Left"
$wsTranslation.Range("E5").Value = "LTR"
$wsTranslation.Range("F5").Value = "RX: <value>"

$wsTranslation.Range("B6").Value = "SingleUseId3"
$wsTranslation.Range("C6").Value = "Default"
$wsTranslation.Range("D6").Value = "Left"
$wsTranslation.Range("E6").Value = "LTR"
$wsTranslation.Range("F6").Value = "default"
